# VM_TestData_Sample1.xlsx edit
# Commit message: "Modified exceptional handling and test data for no code branch"
#
# This adds several new Selenium XPath locator strings used by the test
# automation framework (account search, SSN/FEIN verification, submission
# tab navigation, etc.) into the existing "Xpath" rows of a few sheets, and
# inserts brand-new "Xpath" rows into two sheets that didn't have one yet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# IndexPage (sheet2): add the xpath for the customer/policy search link
# under the existing "ServicePolicy" column (F) on the Xpath row.
# ---------------------------------------------------------------------
$wsIndex = $wb.Worksheets.Item("IndexPage")
$wsIndex.Range("F2").Value = '//a[@ng-click="CustomerAndPolicySearch()"]'

# ---------------------------------------------------------------------
# NewSubmissionPage (sheet3): fill in several previously-empty Xpath
# cells on row 2 with locator strings (some reuse existing locators,
# some are brand new).
# ---------------------------------------------------------------------
$wsNewSub = $wb.Worksheets.Item("NewSubmissionPage")
$wsNewSub.Range("D2").Value = "//*[@id='IdCustomerAccountName']"
$wsNewSub.Range("F2").Value = '//button[@class="btn btn-success pull-right"]'
$wsNewSub.Range("H2").Value = "//*[@id='PrimaryAccountCity']"
$wsNewSub.Range("J2").Value = "//input[@id='PrimaryAccountState']"
$wsNewSub.Range("M2").Value = "//input[@id='IDphysical_aptsuit']"
$wsNewSub.Range("O2").Value = "//*[@id='Id_legal_entity_type']"
$wsNewSub.Range("S2").Value = "//span[contains(text(),'SSN')]"
$wsNewSub.Range("U2").Value = "//span[contains(text(),'FEIN')]"
$wsNewSub.Range("V2").Value = "//*[@id='IdfeinNumber']"
$wsNewSub.Range("AC2").Value = "//*[@ng-click='accountproceed()']"

# ---------------------------------------------------------------------
# ServicePolicy (sheet4): insert a new "Xpath" row (row 2), pushing the
# existing data row down to row 3.
# ---------------------------------------------------------------------
$wsService = $wb.Worksheets.Item("ServicePolicy")
$wsService.Rows.Item(2).Insert()
$wsService.Range("A2").Value = "Xpath"
$wsService.Range("C2").Value = "//input[@ng-model='CustomerSearchCriteria.AccountName']"
$wsService.Range("D2").Value = '//button[@ng-click="customerSearch();"]'
$wsService.Range("E2").Value = "(//a[@class='vam-table-tdArrow'])[2]"

# ---------------------------------------------------------------------
# CustomerInformationPage (sheet5): insert a new "Xpath" row (row 2),
# pushing the existing data row down to row 3.
# ---------------------------------------------------------------------
$wsCustInfo = $wb.Worksheets.Item("CustomerInformationPage")
$wsCustInfo.Rows.Item(2).Insert()
$wsCustInfo.Range("A2").Value = "Xpath"
$wsCustInfo.Range("K2").Value = '//li[@ng-click="setCustomerAccountDetailsTab(''submissions'')"]'
$wsCustInfo.Range("L2").Value = '//a[@ng-click="CreateNewSubmission()"]'

# ---------------------------------------------------------------------
# Restore selections that shifted as a side effect of editing each sheet,
# and make LogInPage the active sheet/tab again.
# ---------------------------------------------------------------------
$wsIndex.Range("F2").Select()
$wsNewSub.Range("AD17").Select()
$wsService.Range("A2").Select()
$wsCustInfo.Range("L2").Select()

$wsLogin = $wb.Worksheets.Item("LogInPage")
$wsLogin.Activate()
$wsLogin.Range("A2").Select()
